$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1817.9131
$ws.Range("I40").Value = 1704.4166
$ws.Range("J40").Value = 1941.7273
$ws.Range("K40").Value = 1704.4166
$ws.Range("L40").Value = 1941.7273
$ws.Range("M40").Value = -1529.4166
$ws.Range("N40").Value = -2291.7273
$ws.Range("H64").Value = 3447
$ws.Range("J64").Value = 3411
$ws.Range("L64").Value = 3411
$ws.Range("N64").Value = -3907
$ws.Range("H67").Value = 3447
$ws.Range("J67").Value = 3411
$ws.Range("L67").Value = 3411
$ws.Range("N67").Value = -5127
$ws.Range("H74").Value = 3564
$ws.Range("I74").Value = 3317.4546
$ws.Range("J74").Value = 3723.5293
$ws.Range("K74").Value = 3317.4546
$ws.Range("L74").Value = 3723.5293
$ws.Range("M74").Value = -2381.4546
$ws.Range("N74").Value = -5595.5293
$ws.Range("H76").Value = 7160.7393
$ws.Range("I76").Value = 5766.5557
$ws.Range("J76").Value = 8057
$ws.Range("K76").Value = 5766.5557
$ws.Range("L76").Value = 8057
$ws.Range("M76").Value = -5451.5557
$ws.Range("N76").Value = -8687
$ws.Range("H77").Value = 3564
$ws.Range("I77").Value = 3317.4546
$ws.Range("J77").Value = 3723.5293
$ws.Range("K77").Value = 16587.273
$ws.Range("L77").Value = 18617.6465
$ws.Range("M77").Value = -11907.273
$ws.Range("N77").Value = -27977.6465
$ws.Range("H79").Value = 7160.7393
$ws.Range("I79").Value = 5766.5557
$ws.Range("J79").Value = 8057
$ws.Range("K79").Value = 5766.5557
$ws.Range("L79").Value = 8057
$ws.Range("M79").Value = -4674.5557
$ws.Range("N79").Value = -10241
$ws.Range("H129").Value = 2608.0303
$ws.Range("I129").Value = 661.3333
$ws.Range("J129").Value = 2802.7
$ws.Range("K129").Value = 1983.9999
$ws.Range("L129").Value = 8408.099999999999
$ws.Range("M129").Value = 3016.0001
$ws.Range("N129").Value = -18408.1
$ws.Range("H132").Value = 3809.3235
$ws.Range("I132").Value = 3697.2334
$ws.Range("J132").Value = 4650
$ws.Range("K132").Value = 11091.7002
$ws.Range("L132").Value = 13950
$ws.Range("M132").Value = -8561.700199999999
$ws.Range("N132").Value = -19010
$ws.Range("H137").Value = 1342
$ws.Range("I137").Value = 977.7222
$ws.Range("J137").Value = 2161.625
$ws.Range("K137").Value = 2933.1666
$ws.Range("L137").Value = 6484.875
$ws.Range("M137").Value = -383.1666
$ws.Range("N137").Value = -11584.875
$ws.Range("H138").Value = 2247.919
$ws.Range("J138").Value = 2728.9473
$ws.Range("L138").Value = 8186.841899999999
$ws.Range("N138").Value = -18466.8419

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5312.4443
$ws.Range("I61").Value = 5476.5
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 5476.5
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -5264.5
$ws.Range("N61").Value = -4424
$ws.Range("H63").Value = 2602.1428
$ws.Range("I63").Value = 2602.1428
$ws.Range("K63").Value = 2602.1428
$ws.Range("M63").Value = -1916.1428
$ws.Range("H66").Value = 2602.1428
$ws.Range("I66").Value = 2602.1428
$ws.Range("K66").Value = 13010.714
$ws.Range("M66").Value = -9578.714
$ws.Range("H132").Value = 22732328
$ws.Range("I132").Value = 31254450
$ws.Range("K132").Value = 93763350
$ws.Range("M132").Value = -93760820
$ws.Range("H136").Value = 5312.4443
$ws.Range("I136").Value = 5476.5
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 16429.5
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -13879.5
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2667.5278
$ws.Range("I134").Value = 1645.4138
$ws.Range("J134").Value = 6902
$ws.Range("K134").Value = 4936.2414
$ws.Range("L134").Value = 20706
$ws.Range("M134").Value = -2401.2414
$ws.Range("N134").Value = -25776

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2642.9375
$ws.Range("J58").Value = 3660
$ws.Range("L58").Value = 3660
$ws.Range("N58").Value = -4066
$ws.Range("H62").Value = 2993.3333
$ws.Range("I62").Value = 2980
$ws.Range("K62").Value = 2980
$ws.Range("M62").Value = -2356
$ws.Range("H65").Value = 2993.3333
$ws.Range("I65").Value = 2980
$ws.Range("K65").Value = 14900
$ws.Range("M65").Value = -11780
$ws.Range("H132").Value = 15983.111
$ws.Range("I132").Value = 20306
$ws.Range("K132").Value = 60918
$ws.Range("M132").Value = -58388
$ws.Range("H134").Value = 2808.625
$ws.Range("I134").Value = 2495.6428
$ws.Range("J134").Value = 4999.5
$ws.Range("K134").Value = 7486.928400000001
$ws.Range("L134").Value = 14998.5
$ws.Range("M134").Value = -4951.928400000001
$ws.Range("N134").Value = -20068.5
$ws.Range("H136").Value = 2642.9375
$ws.Range("J136").Value = 3660
$ws.Range("L136").Value = 10980
$ws.Range("N136").Value = -16080
$ws.Range("H140").Value = 30060
$ws.Range("J140").Value = 30060
$ws.Range("L140").Value = 30060
$ws.Range("N140").Value = -40420

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 96.916664
$ws.Range("I6").Value = 60.272728
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 180.818184
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = -67.818184
$ws.Range("N6").Value = -1726
$ws.Range("H122").Value = 1539.85
$ws.Range("I122").Value = 1276.9166
$ws.Range("J122").Value = 1934.25
$ws.Range("K122").Value = 11492.2494
$ws.Range("L122").Value = 17408.25
$ws.Range("M122").Value = -9042.249400000001
$ws.Range("N122").Value = -22308.25
$ws.Range("H125").Value = 3975
$ws.Range("J125").Value = 4333.3335
$ws.Range("L125").Value = 13000.0005
$ws.Range("N125").Value = -22840.0005
$ws.Range("H134").Value = 3927.6924
$ws.Range("J134").Value = 6655.5557
$ws.Range("L134").Value = 19966.6671
$ws.Range("N134").Value = -30106.6671

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16313.434
$ws.Range("I70").Value = 19648.795
$ws.Range("K70").Value = 19648.795
$ws.Range("M70").Value = -19378.795
$ws.Range("H73").Value = 16313.434
$ws.Range("I73").Value = 19648.795
$ws.Range("K73").Value = 19648.795
$ws.Range("M73").Value = -18712.795
$ws.Range("H80").Value = 12348522
$ws.Range("I80").Value = 17546616
$ws.Range("K80").Value = 17546616
$ws.Range("M80").Value = -17545618
$ws.Range("H83").Value = 12348522
$ws.Range("I83").Value = 17546616
$ws.Range("K83").Value = 87733080
$ws.Range("M83").Value = -87728088
$ws.Range("H132").Value = 6306.381
$ws.Range("I132").Value = 5440.8
$ws.Range("J132").Value = 7093.273
$ws.Range("K132").Value = 16322.4
$ws.Range("L132").Value = 21279.819
$ws.Range("M132").Value = -13792.4
$ws.Range("N132").Value = -26339.819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 17650
$ws.Range("J76").Value = 17650
$ws.Range("L76").Value = 17650
$ws.Range("N76").Value = -18326
$ws.Range("H79").Value = 17650
$ws.Range("J79").Value = 17650
$ws.Range("L79").Value = 17650
$ws.Range("N79").Value = -19990
$ws.Range("H132").Value = 12507925
$ws.Range("I132").Value = 5485.353
$ws.Range("J132").Value = 21748858
$ws.Range("K132").Value = 16456.059
$ws.Range("L132").Value = 65246574
$ws.Range("M132").Value = -13926.059
$ws.Range("N132").Value = -65251634

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1765.5
$ws.Range("I136").Value = 1845.3
$ws.Range("J136").Value = 1632.5
$ws.Range("K136").Value = 5535.9
$ws.Range("L136").Value = 4897.5
$ws.Range("M136").Value = -2985.9
$ws.Range("N136").Value = -9997.5
